$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell's formatting (bold font, border, centered alignment)
# onto the two new header cells before putting values in them, so the new
# columns match the look of the existing header row (B1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for row 2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
